$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the Price column (D) as text so values like "214.27" or "27.820.84" are not
# reinterpreted as numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2: update D,E
$ws.Range("D2").Value = "27.820.84"
$ws.Range("E2").Value = "  +2.52%  "

# Row 3: update D,E
$ws.Range("D3").Value = "1.660.55"
$ws.Range("E3").Value = "  -1.04%  "

# Row 4: update E
$ws.Range("E4").Value = "  +0.25%  "

# Row 5: update D,E
$ws.Range("D5").Value = "214.27"
$ws.Range("E5").Value = "  +0.00%  "

# Row 6: update D,E
$ws.Range("D6").Value = "0.513"
$ws.Range("E6").Value = "  -1.13%  "

# Row 7: update E
$ws.Range("E7").Value = "  +0.26%  "

# Row 8: update D,E
$ws.Range("D8").Value = "23.39"
$ws.Range("E8").Value = "  +2.99%  "

# Row 9: update D,E
$ws.Range("D9").Value = "0.261"
$ws.Range("E9").Value = "  +0.31%  "

# Row 10: update D,E
$ws.Range("D10").Value = "0.0618"
$ws.Range("E10").Value = "  -0.45%  "

# Row 11: update E
$ws.Range("E11").Value = "  -1.07%  "

# Row 12: update D,E
$ws.Range("D12").Value = "1.899.61"
$ws.Range("E12").Value = "  -0.82%  "

# Row 13: update D,E
$ws.Range("D13").Value = "1.660.68"
$ws.Range("E13").Value = "  -1.18%  "

# Row 14: update D,E
$ws.Range("D14").Value = "4.12"
$ws.Range("E14").Value = "  -1.73%  "

# Row 15: update D,E
$ws.Range("D15").Value = "0.548"
$ws.Range("E15").Value = "  -0.21%  "

# Row 16: update D,E
$ws.Range("D16").Value = "65.70"
$ws.Range("E16").Value = "  -1.18%  "

# Row 17: update D,E
$ws.Range("D17").Value = "247.84"
$ws.Range("E17").Value = "  +5.51%  "

# Row 18: update D,E
$ws.Range("D18").Value = "27.825.46"
$ws.Range("E18").Value = "  +2.67%  "

# Row 19: update D,E
$ws.Range("D19").Value = "0.0₃0730"
$ws.Range("E19").Value = "  -1.30%  "

# Row 20: update D,E
$ws.Range("D20").Value = "7.55"
$ws.Range("E20").Value = "  -4.02%  "

# Row 21: update E
$ws.Range("E21").Value = "  +0.23%  "

# Row 22: update D,E
$ws.Range("D22").Value = "4.45"
$ws.Range("E22").Value = "  -1.65%  "

# Row 23: update D,E
$ws.Range("D23").Value = "9.37"
$ws.Range("E23").Value = "  -1.39%  "

# Row 24: update E
$ws.Range("E24").Value = "  -1.49%  "

# Row 25: update D,E
$ws.Range("D25").Value = "146.62"
$ws.Range("E25").Value = "  -0.79%  "

# Row 26: update D,E
$ws.Range("D26").Value = "7.20"
$ws.Range("E26").Value = "  -3.06%  "

# Row 27: update D,E
$ws.Range("D27").Value = "16.15"
$ws.Range("E27").Value = "  -1.08%  "

# Row 28: update E
$ws.Range("E28").Value = "  +0.20%  "

# Row 29: update D,E
$ws.Range("D29").Value = "0.111"
$ws.Range("E29").Value = "  -1.47%  "

# Row 30: update E
$ws.Range("E30").Value = "  +5.95%  "

# Row 31: update D,E
$ws.Range("D31").Value = "0.0498"
$ws.Range("E31").Value = "  -0.49%  "

# Row 32: update E
$ws.Range("E32").Value = "  -0.41%  "

# Row 33: update D,E
$ws.Range("D33").Value = "3.12"
$ws.Range("E33").Value = "  -3.45%  "

# Row 34: update D,E
$ws.Range("D34").Value = "1.405.67"
$ws.Range("E34").Value = "  -8.87%  "

# Row 35: update D,E
$ws.Range("D35").Value = "1.56"
$ws.Range("E35").Value = "  -5.85%  "

# Row 36: update D,E
$ws.Range("D36").Value = "2.40"
$ws.Range("E36").Value = "  +0.24%  "

# Row 37: update D,E
$ws.Range("D37").Value = "0.926"
$ws.Range("E37").Value = "  -1.39%  "

# Row 38: update D,E
$ws.Range("D38").Value = "0.577"
$ws.Range("E38").Value = "  -4.67%  "

# Row 39: update D,E
$ws.Range("D39").Value = "0.0169"
$ws.Range("E39").Value = "  -1.75%  "

# Row 40: update D,E
$ws.Range("D40").Value = "1.03"
$ws.Range("E40").Value = "  -3.33%  "

# Row 41: update D,E
$ws.Range("D41").Value = "69.04"
$ws.Range("E41").Value = "  -0.43%  "

# Row 42: update E
$ws.Range("E42").Value = "  +0.14%  "

# Row 43: update B,C,D,E
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "5.42"
$ws.Range("E43").Value = "  -6.32%  "

# Row 44: update B,C,D,E
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").Value = "2.22"
$ws.Range("E44").Value = "  -1.13%  "

# Row 45: update D,E
$ws.Range("D45").Value = "1.807.45"
$ws.Range("E45").Value = "  -0.81%  "

# Row 46: update D,E
$ws.Range("D46").Value = "0.787"
$ws.Range("E46").Value = "  +1.17%  "

# Row 47: update E
$ws.Range("E47").Value = "  +4.39%  "

# Row 48: update B,C,D,E
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0110"
$ws.Range("E48").Value = "  -0.40%  "

# Row 49: update B,C,D,E
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "88.04"
$ws.Range("E49").Value = "  -2.10%  "

# Row 50: update E
$ws.Range("E50").Value = "  -2.76%  "

# Row 51: update E
$ws.Range("E51").Value = "  -0.35%  "
